$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 295102.3
$ws.Cells.Item(64, 9).Value = 489165.2
$ws.Cells.Item(64, 10).Value = 4008
$ws.Cells.Item(64, 11).Value = 489165.2
$ws.Cells.Item(64, 12).Value = 4008
$ws.Cells.Item(64, 13).Value = -488917.2
$ws.Cells.Item(64, 14).Value = -4504

$ws.Cells.Item(67, 8).Value = 295102.3
$ws.Cells.Item(67, 9).Value = 489165.2
$ws.Cells.Item(67, 10).Value = 4008
$ws.Cells.Item(67, 11).Value = 489165.2
$ws.Cells.Item(67, 12).Value = 4008
$ws.Cells.Item(67, 13).Value = -488307.2
$ws.Cells.Item(67, 14).Value = -5724

$ws.Cells.Item(100, 8).Value = 1328.2941
$ws.Cells.Item(100, 9).Value = 1261
$ws.Cells.Item(100, 10).Value = 1489.8
$ws.Cells.Item(100, 11).Value = 1261
$ws.Cells.Item(100, 12).Value = 1489.8
$ws.Cells.Item(100, 13).Value = -720
$ws.Cells.Item(100, 14).Value = -2571.8

$ws.Cells.Item(137, 8).Value = 2041.4445
$ws.Cells.Item(137, 9).Value = 1377.3158
$ws.Cells.Item(137, 11).Value = 4131.9474
$ws.Cells.Item(137, 13).Value = -1581.9474

$ws.Cells.Item(141, 8).Value = 4074.5715
$ws.Cells.Item(141, 9).Value = 1590.9783
$ws.Cells.Item(141, 10).Value = 10794.883
$ws.Cells.Item(141, 11).Value = 4772.9349
$ws.Cells.Item(141, 12).Value = 32384.649
$ws.Cells.Item(141, 13).Value = 407.0650999999998
$ws.Cells.Item(141, 14).Value = -42744.649

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 5001277.5
$ws.Cells.Item(122, 9).Value = 1260.8422
$ws.Cells.Item(122, 10).Value = 20834664
$ws.Cells.Item(122, 11).Value = 3782.5266
$ws.Cells.Item(122, 12).Value = 62503992
$ws.Cells.Item(122, 13).Value = -1332.5266
$ws.Cells.Item(122, 14).Value = -62508892

$ws.Cells.Item(132, 8).Value = 5674.2
$ws.Cells.Item(132, 9).Value = 2350.3333
$ws.Cells.Item(132, 10).Value = 10660
$ws.Cells.Item(132, 11).Value = 7050.999899999999
$ws.Cells.Item(132, 12).Value = 31980
$ws.Cells.Item(132, 13).Value = -4520.999899999999
$ws.Cells.Item(132, 14).Value = -37040

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 26894.637
$ws.Cells.Item(31, 9).Value = 33515.117
$ws.Cells.Item(31, 10).Value = 4385
$ws.Cells.Item(31, 11).Value = 33515.117
$ws.Cells.Item(31, 12).Value = 4385
$ws.Cells.Item(31, 13).Value = -33220.117
$ws.Cells.Item(31, 14).Value = -4975

$ws.Cells.Item(34, 8).Value = 26894.637
$ws.Cells.Item(34, 9).Value = 33515.117
$ws.Cells.Item(34, 10).Value = 4385
$ws.Cells.Item(34, 11).Value = 33515.117
$ws.Cells.Item(34, 12).Value = 4385
$ws.Cells.Item(34, 13).Value = -33313.117
$ws.Cells.Item(34, 14).Value = -4789

$ws.Cells.Item(58, 8).Value = 1655207.2
$ws.Cells.Item(58, 9).Value = 2675443.2
$ws.Cells.Item(58, 10).Value = 3396.6667
$ws.Cells.Item(58, 11).Value = 2675443.2
$ws.Cells.Item(58, 12).Value = 3396.6667
$ws.Cells.Item(58, 13).Value = -2675240.2
$ws.Cells.Item(58, 14).Value = -3802.6667

$ws.Cells.Item(92, 8).Value = 49800
$ws.Cells.Item(92, 10).Value = 49800
$ws.Cells.Item(92, 12).Value = 49800
$ws.Cells.Item(92, 14).Value = -54792

$ws.Cells.Item(122, 8).Value = 9392.5
$ws.Cells.Item(122, 9).Value = 4706.4165
$ws.Cells.Item(122, 10).Value = 15015.8
$ws.Cells.Item(122, 11).Value = 14119.2495
$ws.Cells.Item(122, 12).Value = 45047.39999999999
$ws.Cells.Item(122, 13).Value = -11669.2495
$ws.Cells.Item(122, 14).Value = -49947.39999999999

$ws.Cells.Item(136, 8).Value = 1655207.2
$ws.Cells.Item(136, 9).Value = 2675443.2
$ws.Cells.Item(136, 10).Value = 3396.6667
$ws.Cells.Item(136, 11).Value = 8026329.600000001
$ws.Cells.Item(136, 12).Value = 10190.0001
$ws.Cells.Item(136, 13).Value = -8023779.600000001
$ws.Cells.Item(136, 14).Value = -15290.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 4276604.5
$ws.Cells.Item(5, 9).Value = 343.2963
$ws.Cells.Item(5, 10).Value = 13898192
$ws.Cells.Item(5, 11).Value = 1029.8889
$ws.Cells.Item(5, 12).Value = 41694576
$ws.Cells.Item(5, 13).Value = -917.8888999999999
$ws.Cells.Item(5, 14).Value = -41694800

$ws.Cells.Item(68, 8).Value = 5342.609
$ws.Cells.Item(68, 9).Value = 675.8
$ws.Cells.Item(68, 10).Value = 6638.9443
$ws.Cells.Item(68, 11).Value = 2027.4
$ws.Cells.Item(68, 12).Value = 19916.8329
$ws.Cells.Item(68, 13).Value = -1216.4
$ws.Cells.Item(68, 14).Value = -21538.8329

$ws.Cells.Item(71, 8).Value = 5342.609
$ws.Cells.Item(71, 9).Value = 675.8
$ws.Cells.Item(71, 10).Value = 6638.9443
$ws.Cells.Item(71, 11).Value = 6082.2
$ws.Cells.Item(71, 12).Value = 59750.4987
$ws.Cells.Item(71, 13).Value = -2026.2
$ws.Cells.Item(71, 14).Value = -67862.4987

$ws.Cells.Item(92, 8).Value = 855.7778
$ws.Cells.Item(92, 9).Value = 797.3333
$ws.Cells.Item(92, 10).Value = 885
$ws.Cells.Item(92, 11).Value = 2391.9999
$ws.Cells.Item(92, 12).Value = 2655
$ws.Cells.Item(92, 13).Value = -1143.9999
$ws.Cells.Item(92, 14).Value = -5151

$ws.Cells.Item(113, 8).Value = 709.1316
$ws.Cells.Item(113, 9).Value = 705.55554
$ws.Cells.Item(113, 10).Value = 726.46155
$ws.Cells.Item(113, 11).Value = 2116.66662
$ws.Cells.Item(113, 12).Value = 2179.38465
$ws.Cells.Item(113, 13).Value = 53.33338000000003
$ws.Cells.Item(113, 14).Value = -6519.38465

$ws.Cells.Item(135, 8).Value = 4276604.5
$ws.Cells.Item(135, 9).Value = 343.2963
$ws.Cells.Item(135, 10).Value = 13898192
$ws.Cells.Item(135, 11).Value = 3089.6667
$ws.Cells.Item(135, 12).Value = 125083728
$ws.Cells.Item(135, 13).Value = -554.6666999999998
$ws.Cells.Item(135, 14).Value = -125088798

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 11900
$ws.Cells.Item(92, 10).Value = 11900
$ws.Cells.Item(92, 12).Value = 11900
$ws.Cells.Item(92, 14).Value = -15644

$ws.Cells.Item(97, 8).Value = 1230.5
$ws.Cells.Item(97, 9).Value = 831.44446
$ws.Cells.Item(97, 10).Value = 2210
$ws.Cells.Item(97, 11).Value = 831.44446
$ws.Cells.Item(97, 12).Value = 2210
$ws.Cells.Item(97, 13).Value = -335.44446
$ws.Cells.Item(97, 14).Value = -3202

$ws.Cells.Item(102, 8).Value = 5102.161
$ws.Cells.Item(102, 9).Value = 4669.9546
$ws.Cells.Item(102, 10).Value = 6158.6665
$ws.Cells.Item(102, 11).Value = 4669.9546
$ws.Cells.Item(102, 12).Value = 6158.6665
$ws.Cells.Item(102, 13).Value = -3047.9546
$ws.Cells.Item(102, 14).Value = -9402.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9225.625
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 9225.625
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 9225.625
$ws.Cells.Item(7, 13).Value = $null
$ws.Cells.Item(7, 14).Value = -9449.625

$ws.Cells.Item(55, 8).Value = 259.1579
$ws.Cells.Item(55, 9).Value = 227
$ws.Cells.Item(55, 10).Value = 303.375
$ws.Cells.Item(55, 11).Value = 227
$ws.Cells.Item(55, 12).Value = 303.375
$ws.Cells.Item(55, 13).Value = -54
$ws.Cells.Item(55, 14).Value = -649.375

$ws.Cells.Item(93, 8).Value = 453.69232
$ws.Cells.Item(93, 9).Value = 437.6875
$ws.Cells.Item(93, 10).Value = 479.3
$ws.Cells.Item(93, 11).Value = 437.6875
$ws.Cells.Item(93, 12).Value = 479.3
$ws.Cells.Item(93, 13).Value = 810.3125
$ws.Cells.Item(93, 14).Value = -2975.3

$ws.Cells.Item(122, 8).Value = 7600.2764
$ws.Cells.Item(122, 9).Value = 7228
$ws.Cells.Item(122, 10).Value = 8818.637000000001
$ws.Cells.Item(122, 11).Value = 21684
$ws.Cells.Item(122, 12).Value = 26455.911
$ws.Cells.Item(122, 13).Value = -19234
$ws.Cells.Item(122, 14).Value = -31355.911

$ws.Cells.Item(126, 8).Value = 9225.625
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 9225.625
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 27676.875
$ws.Cells.Item(126, 13).Value = $null
$ws.Cells.Item(126, 14).Value = -32616.875

$ws.Cells.Item(132, 8).Value = 5971.0938
$ws.Cells.Item(132, 9).Value = 6799.85
$ws.Cells.Item(132, 10).Value = 4589.8335
$ws.Cells.Item(132, 11).Value = 20399.55
$ws.Cells.Item(132, 12).Value = 13769.5005
$ws.Cells.Item(132, 13).Value = -17869.55
$ws.Cells.Item(132, 14).Value = -18829.5005

$ws.Cells.Item(136, 8).Value = 2978.5542
$ws.Cells.Item(136, 9).Value = 1745.5834
$ws.Cells.Item(136, 10).Value = 6195
$ws.Cells.Item(136, 11).Value = 5236.7502
$ws.Cells.Item(136, 12).Value = 18585
$ws.Cells.Item(136, 13).Value = -2686.7502
$ws.Cells.Item(136, 14).Value = -23685

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1575
$ws.Cells.Item(96, 9).Value = 1599
$ws.Cells.Item(96, 11).Value = 1599
$ws.Cells.Item(96, 13).Value = -226

$ws.Cells.Item(122, 8).Value = 2885.8276
$ws.Cells.Item(122, 9).Value = 2315.36
$ws.Cells.Item(122, 11).Value = 6946.08
$ws.Cells.Item(122, 13).Value = -4496.08

$ws.Cells.Item(132, 8).Value = 1285.1666
$ws.Cells.Item(132, 9).Value = 493.16666
$ws.Cells.Item(132, 10).Value = 3133.1667
$ws.Cells.Item(132, 11).Value = 1479.49998
$ws.Cells.Item(132, 12).Value = 9399.500100000001
$ws.Cells.Item(132, 13).Value = 1050.50002
$ws.Cells.Item(132, 14).Value = -14459.5001

$ws.Cells.Item(136, 8).Value = 3726.875
$ws.Cells.Item(136, 9).Value = 2646.2812
$ws.Cells.Item(136, 10).Value = 5888.0625
$ws.Cells.Item(136, 11).Value = 7938.8436
$ws.Cells.Item(136, 12).Value = 17664.1875
$ws.Cells.Item(136, 13).Value = -5388.8436
$ws.Cells.Item(136, 14).Value = -22764.1875
